# Generate Report for Handoff
# Updates the localization-status report: the dbca198b-...-md file moved
# from "Handed back: in sync with en-US" to "Ready for handoff" with a new
# handoff generated (and a detail note that the existing handback file is
# stale relative to the newly generated handoff).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/5ab631ea6f0732e30e27ee0ec371ea7e83c88c6d/e2e/dbca198b-97a7-4d9e-9af7-45a3899e2554.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/de79c7a70e273d13ff9b4d3a63b96482fc2b8f51/e2e/dbca198b-97a7-4d9e-9af7-45a3899e2554.md."

# --- Overview sheet: row 3 is the dbca198b-...-md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-13 04:58:55"

# --- zh-cn sheet: row 3 is the dbca198b-...-md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-13 04:58:47"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1

# --- de-de sheet: row 3 is the dbca198b-...-md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-13 04:58:55"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1

Write-Output "Report updated for handoff"
